$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the Price/Volume columns as Text before writing the new
# values. This stops Excel from auto-converting numeric-looking strings such
# as "1.037" or "0.06700" into floating point numbers (which would silently
# drop significant trailing zeros and change the cell type). The style is
# reset back to Normal afterwards so the cells keep their original formatting.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.931.11'
$ws.Range("E2").Value = '  +3.36%  '

$ws.Range("D3").Value = '1.918.80'
$ws.Range("E3").Value = '  +3.20%  '

$ws.Range("D4").Value = '1.037'
$ws.Range("E4").Value = '  +3.31%  '

$ws.Range("D5").Value = '323.02'
$ws.Range("E5").Value = '  +3.40%  '

$ws.Range("D6").Value = '1.037'
$ws.Range("E6").Value = '  +3.39%  '

$ws.Range("D7").Value = '0.5237'
$ws.Range("E7").Value = '  +1.81%  '

$ws.Range("D8").Value = '0.3986'
$ws.Range("E8").Value = '  +3.98%  '

$ws.Range("D9").Value = '0.08495'
$ws.Range("E9").Value = '  +3.05%  '

$ws.Range("D10").Value = '1.146'
$ws.Range("E10").Value = '  +3.31%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '42.78'
$ws.Range("E11").Value = '  +3.18%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.357'
$ws.Range("E12").Value = '  +2.69%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.922.10'
$ws.Range("E13").Value = '  +3.06%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '20.81'
$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.393'
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.043'
$ws.Range("E16").Value = '  +3.87%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001126'
$ws.Range("E17").Value = '  +2.70%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '92.15'
$ws.Range("E18").Value = '  +1.76%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.06807'
$ws.Range("E19").Value = '  +2.32%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '18.15'
$ws.Range("E20").Value = '  +2.73%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.035'
$ws.Range("E21").Value = '  +3.20%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.141'
$ws.Range("E22").Value = '  +2.31%  '

$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '28.928.24'
$ws.Range("E23").Value = '  +3.23%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  +2.67%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.336'
$ws.Range("E25").Value = '  +3.71%  '

$ws.Range("D26").Value = '2.132.31'
$ws.Range("E26").Value = '  +2.74%  '

$ws.Range("D27").Value = '164.18'
$ws.Range("E27").Value = '  +3.80%  '

$ws.Range("D28").Value = '21.19'
$ws.Range("E28").Value = '  +3.60%  '

$ws.Range("D29").Value = '2.491'
$ws.Range("E29").Value = '  -1.01%  '

$ws.Range("D30").Value = '128.55'
$ws.Range("E30").Value = '  +3.21%  '

$ws.Range("D31").Value = '0.1067'
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").Value = '1.061'
$ws.Range("E32").Value = '  +3.22%  '

$ws.Range("D33").Value = '6.032'
$ws.Range("E33").Value = '  +1.02%  '

$ws.Range("D34").Value = '3.707'
$ws.Range("E34").Value = '  +2.98%  '

$ws.Range("D35").Value = '9.591'
$ws.Range("E35").Value = '  +2.62%  '

$ws.Range("D36").Value = '0.02499'
$ws.Range("E36").Value = '  +3.48%  '

$ws.Range("D37").Value = '0.06700'
$ws.Range("E37").Value = '  +3.13%  '

$ws.Range("D38").Value = '0.2254'
$ws.Range("E38").Value = '  +3.96%  '

$ws.Range("D39").Value = '0.6627'
$ws.Range("E39").Value = '  +1.25%  '

$ws.Range("D40").Value = '1.272'
$ws.Range("E40").Value = '  +3.81%  '

$ws.Range("D41").Value = '1.208'
$ws.Range("E41").Value = '  +1.11%  '

$ws.Range("D42").Value = '5.085'
$ws.Range("E42").Value = '  +1.56%  '

$ws.Range("D43").Value = '11.33'
$ws.Range("E43").Value = '  +1.73%  '

$ws.Range("D44").Value = '0.6227'
$ws.Range("E44").Value = '  +1.48%  '

$ws.Range("D45").Value = '13.33'
$ws.Range("E45").Value = '  +2.84%  '

$ws.Range("D46").Value = '3.785'
$ws.Range("E46").Value = '  +3.43%  '

$ws.Range("D47").Value = '1.314'
$ws.Range("E47").Value = '  +2.66%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.046'
$ws.Range("E48").Value = '  +1.96%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.257'
$ws.Range("E49").Value = '  +3.53%  '

$ws.Range("D50").Value = '123.39'
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("D51").Value = '0.07006'
$ws.Range("E51").Value = '  +2.38%  '

# Restore the default cell style now that the text values are in place.
$dataRange.Style = "Normal"
